$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 73, shifting existing rows 73-77 down to 75-79.
$ws.Rows.Item(73).Resize(2).Insert()

# Populate the new row 73 (Primera quality, newest report date 2021-11-05 = 44505)
$ws.Cells.Item(73, 1).Value = 11
$ws.Cells.Item(73, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(73, 3).Value = "Bíobío"
$ws.Cells.Item(73, 4).Value = 44505
$ws.Cells.Item(73, 4).NumberFormat = $ws.Cells.Item(75, 4).NumberFormat
$ws.Cells.Item(73, 5).Value = 8
$ws.Cells.Item(73, 6).Value = "Fruta"
$ws.Cells.Item(73, 7).Value = 100108
$ws.Cells.Item(73, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(73, 9).Value = 100108002
$ws.Cells.Item(73, 10).Value = "Mango"
$ws.Cells.Item(73, 11).Value = "Sin especificar"
$ws.Cells.Item(73, 12).Value = "Primera"
$ws.Cells.Item(73, 13).Value = 200
$ws.Cells.Item(73, 14).Value = 7000
$ws.Cells.Item(73, 15).Value = 7500
$ws.Cells.Item(73, 16).Value = 7250
$ws.Cells.Item(73, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(73, 18).Value = "Perú"
$ws.Cells.Item(73, 19).Value = 1812
$ws.Cells.Item(73, 20).Value = 4

# Populate the new row 74 (Segunda quality, same newest report date)
$ws.Cells.Item(74, 1).Value = 11
$ws.Cells.Item(74, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(74, 3).Value = "Bíobío"
$ws.Cells.Item(74, 4).Value = 44505
$ws.Cells.Item(74, 4).NumberFormat = $ws.Cells.Item(75, 4).NumberFormat
$ws.Cells.Item(74, 5).Value = 8
$ws.Cells.Item(74, 6).Value = "Fruta"
$ws.Cells.Item(74, 7).Value = 100108
$ws.Cells.Item(74, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(74, 9).Value = 100108002
$ws.Cells.Item(74, 10).Value = "Mango"
$ws.Cells.Item(74, 11).Value = "Sin especificar"
$ws.Cells.Item(74, 12).Value = "Segunda"
$ws.Cells.Item(74, 13).Value = 100
$ws.Cells.Item(74, 14).Value = 6500
$ws.Cells.Item(74, 15).Value = 6500
$ws.Cells.Item(74, 16).Value = 6500
$ws.Cells.Item(74, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(74, 18).Value = "Perú"
$ws.Cells.Item(74, 19).Value = 1625
$ws.Cells.Item(74, 20).Value = 4
